$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D hold text that looks numeric (e.g. "1.002", "286.70").
# Excel auto-converts such text to a real number on assignment unless the
# cell is pre-formatted as Text ("@"). We flip the format to Text, assign
# the literal string, then restore the "Normal" cell style so no stray
# number-format style sticks to the cell (keeps styles.xml/cell style index
# identical to the original).

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

$ws.Range("D2").Value = "21.851.64"
$ws.Range("E2").Value = "  +6.43%  "

$ws.Range("D3").Value = "1.573.65"
$ws.Range("E3").Value = "  +6.98%  "

$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  -0.53%  "

$ws.Range("D5").Value = "0.9883"
$ws.Range("E5").Value = "  +2.50%  "

$ws.Range("D6").Value = "286.70"
$ws.Range("E6").Value = "  +3.75%  "

$ws.Range("D7").Value = "0.3695"
$ws.Range("E7").Value = "  +1.25%  "

$ws.Range("D8").Value = "0.3281"
$ws.Range("E8").Value = "  +7.28%  "

$ws.Range("D9").Value = "1.142"
$ws.Range("E9").Value = "  +7.83%  "

$ws.Range("D10").Value = "41.61"
$ws.Range("E10").Value = "  +3.59%  "

$ws.Range("D11").Value = "0.07036"
$ws.Range("E11").Value = "  +6.12%  "

$ws.Range("D12").Value = "0.9984"
$ws.Range("E12").Value = "  -0.24%  "

$ws.Range("D13").Value = "20.06"
$ws.Range("E13").Value = "  +10.62%  "

$ws.Range("D14").Value = "5.842"
$ws.Range("E14").Value = "  +6.76%  "

$ws.Range("D15").Value = "6.524"
$ws.Range("E15").Value = "  +5.58%  "

$ws.Range("D16").Value = "0.00001072"
$ws.Range("E16").Value = "  +4.06%  "

$ws.Range("D17").Value = "0.9880"
$ws.Range("E17").Value = "  +2.58%  "

$ws.Range("D18").Value = "1.570.79"
$ws.Range("E18").Value = "  +6.56%  "

$ws.Range("D19").Value = "0.06326"
$ws.Range("E19").Value = "  +6.84%  "

$ws.Range("D20").Value = "75.36"
$ws.Range("E20").Value = "  +8.86%  "

$ws.Range("D21").Value = "16.14"
$ws.Range("E21").Value = "  +11.13%  "

$ws.Range("D22").Value = "5.862"
$ws.Range("E22").Value = "  +7.36%  "

$ws.Range("D23").Value = "11.62"
$ws.Range("E23").Value = "  +5.01%  "

$ws.Range("B24").Value = "WrappedBTC"
$ws.Range("C24").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D24").Value = "21.881.17"
$ws.Range("E24").Value = "  +6.28%  "

$ws.Range("B25").Value = "Toncoin"
$ws.Range("C25").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D25").Value = "2.357"
$ws.Range("E25").Value = "  +5.02%  "

$ws.Range("D26").Value = "2.402"
$ws.Range("E26").Value = "  +12.66%  "

$ws.Range("E27").Value = "  +6.89%  "

$ws.Range("D28").Value = "18.52"
$ws.Range("E28").Value = "  +7.44%  "

$ws.Range("D29").Value = "1.748.27"
$ws.Range("E29").Value = "  +7.26%  "

$ws.Range("D30").Value = "120.27"
$ws.Range("E30").Value = "  +5.41%  "

$ws.Range("D31").Value = "4.141"
$ws.Range("E31").Value = "  +4.82%  "

$ws.Range("D32").Value = "0.9178"
$ws.Range("E32").Value = "  +13.01%  "

$ws.Range("D33").Value = "5.453"
$ws.Range("E33").Value = "  +9.79%  "

$ws.Range("D34").Value = "0.08204"
$ws.Range("E34").Value = "  +3.26%  "

$ws.Range("D35").Value = "1.611"
$ws.Range("E35").Value = "  +4.69%  "

$ws.Range("B36").Value = "FraxShare"
$ws.Range("C36").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D36").Value = "8.688"
$ws.Range("E36").Value = "  +13.53%  "

$ws.Range("B37").Value = "InternetComputer(DFINITY)"
$ws.Range("C37").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D37").Value = "5.135"
$ws.Range("E37").Value = "  +8.84%  "

$ws.Range("D38").Value = "0.06119"
$ws.Range("E38").Value = "  +5.27%  "

$ws.Range("B39").Value = "Aptos"
$ws.Range("C39").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D39").Value = "11.65"
$ws.Range("E39").Value = "  +11.60%  "

$ws.Range("B40").Value = "TrustWalletToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D40").Value = "1.231"
$ws.Range("E40").Value = "  +1.90%  "

$ws.Range("D41").Value = "0.02181"
$ws.Range("E41").Value = "  +7.10%  "

$ws.Range("D42").Value = "0.2005"
$ws.Range("E42").Value = "  +6.64%  "

$ws.Range("D43").Value = "0.9879"
$ws.Range("E43").Value = "  +2.56%  "

$ws.Range("D44").Value = "0.5781"
$ws.Range("E44").Value = "  +9.40%  "

$ws.Range("D45").Value = "12.92"
$ws.Range("E45").Value = "  +6.44%  "

$ws.Range("D46").Value = "3.634"
$ws.Range("E46").Value = "  +3.60%  "

$ws.Range("D47").Value = "0.5646"
$ws.Range("E47").Value = "  +8.74%  "

$ws.Range("D48").Value = "125.15"
$ws.Range("E48").Value = "  +6.02%  "

$ws.Range("E49").Value = "  +7.17%  "

$ws.Range("D50").Value = "0.06764"
$ws.Range("E50").Value = "  +4.65%  "

$ws.Range("D51").Value = "72.07"
$ws.Range("E51").Value = "  +7.49%  "

# Restore the default "Normal" style on the touched D cells so their style
# index matches the original (unstyled) cells.
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").Style = "Normal"
$ws.Range("D4").Style = "Normal"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").Style = "Normal"
$ws.Range("D7").Style = "Normal"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").Style = "Normal"
$ws.Range("D18").Style = "Normal"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").Style = "Normal"
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").Style = "Normal"
$ws.Range("D38").Style = "Normal"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").Style = "Normal"
$ws.Range("D48").Style = "Normal"
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").Style = "Normal"
